$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only touch cells whose content actually changes (or new cells), and in the
# exact order needed so newly-introduced shared strings land in the same
# sequence as the target workbook. Cells that keep their original text
# (header row A1:E1, B4, C2, E4's text) are intentionally left untouched.

# Row 3 (rtx) - text changes
$ws.Range("A3").Value = "rtx 4060 ti 16 gb"
$ws.Range("B3").Value = "3060 3070 4070 3080 8"

# Row 5 (new - ryzen)
$ws.Range("A5").Value = "ryzen 5800x3d"

# Row 6 (new - logitech)
$ws.Range("A6").Value = "logitech g733"

# Row 2 (iphone) - text changes
$ws.Range("A2").Value = "iphone 15 pro max 256 gb"
$ws.Range("B2").Value = "mini watch 11 12 13 14 usado"

# Row 4 (macbook) - text change
$ws.Range("A4").Value = "macbook 16 gb m2"

# Row 5 / Row 6 banned-words continued
$ws.Range("B5").Value = "5700x"
$ws.Range("B6").Value = "g535"

# Banned websites column
$ws.Range("C3").Value = "patoloco aliexpress techinn shopee"
$ws.Range("C4").Value = "aliexpress tiendamia shopee enjoei ebay"
$ws.Range("C5").Value = "aliexpress tiendamia shopee enjoei ebay"
$ws.Range("C6").Value = "aliexpress tiendamia shopee enjoei ebay"

# Numeric values (do not affect shared strings / their order)
$ws.Range("D2").Value = 8000
$ws.Range("E2").Value = 9500
$ws.Range("D3").Value = 2500
$ws.Range("E3").Value = 3500
$ws.Range("D4").Value = 5500
$ws.Range("D5").Value = 1500
$ws.Range("E5").Value = 2500
$ws.Range("D6").Value = 600
$ws.Range("E6").Value = 920

# Underline formatting on specific cells
$ws.Range("E4").Font.Underline = $true
$ws.Range("D6").Font.Underline = $true
$ws.Range("E6").Font.Underline = $true

# Column C width (stored width ends up at 110)
$ws.Columns("C").ColumnWidth = 109.15

# Selection
$ws.Range("C7").Select()
